$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Mmrn2"
$ws.Range("C2").Value = "Clec14a"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 41.006619
$ws.Range("H2").Value = 123.019857
$ws.Range("I2").Value = 0.9349081063755518
$ws.Range("J2").Value = 0.9349081063755517
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 9.929107999999999
$ws.Range("N2").Value = 29.787324
$ws.Range("O2").Value = 0.8066781947969484
$ws.Range("P2").Value = 0.8066781947969485
$ws.Range("Q2").Value = 407.159148765852
$ws.Range("R2").Value = 3664.432338892668
$ws.Range("S2").Value = 0.7541699835520634
$ws.Range("T2").Value = 0.7541699835520634

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Mmrn2"
$ws.Range("C3").Value = "Clec14a"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 41.006619
$ws.Range("H3").Value = 123.019857
$ws.Range("I3").Value = 0.9349081063755518
$ws.Range("J3").Value = 0.9349081063755517
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.894331
$ws.Range("N3").Value = 5.682993
$ws.Range("O3").Value = 0.1539025974365369
$ws.Range("P3").Value = 0.1539025974365369
$ws.Range("Q3").Value = 77.680109576889
$ws.Range("R3").Value = 699.1209861920009
$ws.Range("S3").Value = 0.1438847859356716
$ws.Range("T3").Value = 0.1438847859356716

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Mmrn2"
$ws.Range("C4").Value = "Clec14a"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 41.006619
$ws.Range("H4").Value = 123.019857
$ws.Range("I4").Value = 0.9349081063755518
$ws.Range("J4").Value = 0.9349081063755517
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.01185166666666667
$ws.Range("N4").Value = 0.035555
$ws.Range("O4").Value = 0.000962874114371788
$ws.Range("P4").Value = 0.000962874114371788
$ws.Range("Q4").Value = 0.4859967795150001
$ws.Range("R4").Value = 4.373971015635001
$ws.Range("S4").Value = 0.0009001988149453648
$ws.Range("T4").Value = 0.0009001988149453647

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Mmrn2"
$ws.Range("C5").Value = "Clec14a"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 41.006619
$ws.Range("H5").Value = 123.019857
$ws.Range("I5").Value = 0.9349081063755518
$ws.Range("J5").Value = 0.9349081063755517
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.4733450000000001
$ws.Range("N5").Value = 1.420035
$ws.Range("O5").Value = 0.03845633365214293
$ws.Range("P5").Value = 0.03845633365214293
$ws.Range("Q5").Value = 19.410278070555
$ws.Range("R5").Value = 174.692502634995
$ws.Range("S5").Value = 0.03595313807287136
$ws.Range("T5").Value = 0.03595313807287135

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Mmrn2"
$ws.Range("C6").Value = "Clec14a"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 2.068862333333333
$ws.Range("H6").Value = 6.206586999999999
$ws.Range("I6").Value = 0.0471679015138598
$ws.Range("J6").Value = 0.04716790151385979
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 9.929107999999999
$ws.Range("N6").Value = 29.787324
$ws.Range("O6").Value = 0.8066781947969484
$ws.Range("P6").Value = 0.8066781947969485
$ws.Range("Q6").Value = 20.54195754479866
$ws.Range("R6").Value = 184.877617903188
$ws.Range("S6").Value = 0.03804931764556067
$ws.Range("T6").Value = 0.03804931764556067

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Mmrn2"
$ws.Range("C7").Value = "Clec14a"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 2.068862333333333
$ws.Range("H7").Value = 6.206586999999999
$ws.Range("I7").Value = 0.0471679015138598
$ws.Range("J7").Value = 0.04716790151385979
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.894331
$ws.Range("N7").Value = 5.682993
$ws.Range("O7").Value = 0.1539025974365369
$ws.Range("P7").Value = 0.1539025974365369
$ws.Range("Q7").Value = 3.919110052765666
$ws.Range("R7").Value = 35.27199047489099
$ws.Range("S7").Value = 0.007259262558613785
$ws.Range("T7").Value = 0.007259262558613783

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Mmrn2"
$ws.Range("C8").Value = "Clec14a"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 2.068862333333333
$ws.Range("H8").Value = 6.206586999999999
$ws.Range("I8").Value = 0.0471679015138598
$ws.Range("J8").Value = 0.04716790151385979
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.01185166666666667
$ws.Range("N8").Value = 0.035555
$ws.Range("O8").Value = 0.000962874114371788
$ws.Range("P8").Value = 0.000962874114371788
$ws.Range("Q8").Value = 0.02451946675388889
$ws.Range("R8").Value = 0.220675200785
$ws.Range("S8").Value = 0.00004541675139693348
$ws.Range("T8").Value = 0.00004541675139693346

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Mmrn2"
$ws.Range("C9").Value = "Clec14a"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 2.068862333333333
$ws.Range("H9").Value = 6.206586999999999
$ws.Range("I9").Value = 0.0471679015138598
$ws.Range("J9").Value = 0.04716790151385979
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.4733450000000001
$ws.Range("N9").Value = 1.420035
$ws.Range("O9").Value = 0.03845633365214293
$ws.Range("P9").Value = 0.03845633365214293
$ws.Range("Q9").Value = 0.9792856411716667
$ws.Range("R9").Value = 8.813570770544999
$ws.Range("S9").Value = 0.00181390455828841
$ws.Range("T9").Value = 0.00181390455828841

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Mmrn2"
$ws.Range("C10").Value = "Clec14a"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.7861760000000001
$ws.Range("H10").Value = 2.358528
$ws.Range("I10").Value = 0.01792399211058844
$ws.Range("J10").Value = 0.01792399211058843
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 9.929107999999999
$ws.Range("N10").Value = 29.787324
$ws.Range("O10").Value = 0.8066781947969484
$ws.Range("P10").Value = 0.8066781947969485
$ws.Range("Q10").Value = 7.806026411008
$ws.Range("R10").Value = 70.25423769907201
$ws.Range("S10").Value = 0.01445889359932422
$ws.Range("T10").Value = 0.01445889359932422

# Row 11
$ws.Range("A11").Value = "sCs"
$ws.Range("B11").Value = "Mmrn2"
$ws.Range("C11").Value = "Clec14a"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.7861760000000001
$ws.Range("H11").Value = 2.358528
$ws.Range("I11").Value = 0.01792399211058844
$ws.Range("J11").Value = 0.01792399211058843
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 1.894331
$ws.Range("N11").Value = 5.682993
$ws.Range("O11").Value = 0.1539025974365369
$ws.Range("P11").Value = 0.1539025974365369
$ws.Range("Q11").Value = 1.489277568256
$ws.Range("R11").Value = 13.403498114304
$ws.Range("S11").Value = 0.002758548942251556
$ws.Range("T11").Value = 0.002758548942251555

# Row 12
$ws.Range("A12").Value = "sCs"
$ws.Range("B12").Value = "Mmrn2"
$ws.Range("C12").Value = "Clec14a"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.7861760000000001
$ws.Range("H12").Value = 2.358528
$ws.Range("I12").Value = 0.01792399211058844
$ws.Range("J12").Value = 0.01792399211058843
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.01185166666666667
$ws.Range("N12").Value = 0.035555
$ws.Range("O12").Value = 0.000962874114371788
$ws.Range("P12").Value = 0.000962874114371788
$ws.Range("Q12").Value = 0.009317495893333336
$ws.Range("R12").Value = 0.08385746304000001
$ws.Range("S12").Value = 0.00001725854802948976
$ws.Range("T12").Value = 0.00001725854802948976

# Row 13
$ws.Range("A13").Value = "sCs"
$ws.Range("B13").Value = "Mmrn2"
$ws.Range("C13").Value = "Clec14a"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.7861760000000001
$ws.Range("H13").Value = 2.358528
$ws.Range("I13").Value = 0.01792399211058844
$ws.Range("J13").Value = 0.01792399211058843
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.4733450000000001
$ws.Range("N13").Value = 1.420035
$ws.Range("O13").Value = 0.03845633365214293
$ws.Range("P13").Value = 0.03845633365214293
$ws.Range("Q13").Value = 0.3721324787200001
$ws.Range("R13").Value = 3.349192308480001
$ws.Range("S13").Value = 0.0006892910209831666
$ws.Range("T13").Value = 0.0006892910209831665
